# Insert a new data row for row 23 (a new weekly price record), pushing the
# existing rows 23-93 down to 24-94 (dimension grows from A1:R93 to A1:R94).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new record's data. The
# non-date/value columns (market id/name/region, category info, units,
# origin, classification) are identical across all rows in this block, so
# copy them down from the row directly below (row 24), which still holds the
# original row-23 data after the shift.
$ws.Range("A23").Value = $ws.Range("A24").Value2
$ws.Range("B23").Value = $ws.Range("B24").Value2
$ws.Range("C23").Value = $ws.Range("C24").Value2
$ws.Range("D23").Value = 44575
$ws.Range("E23").Value = $ws.Range("E24").Value2
$ws.Range("F23").Value = $ws.Range("F24").Value2
$ws.Range("G23").Value = $ws.Range("G24").Value2
$ws.Range("H23").Value = $ws.Range("H24").Value2
$ws.Range("I23").Value = $ws.Range("I24").Value2
$ws.Range("J23").Value = 400
$ws.Range("K23").Value = 22000
$ws.Range("L23").Value = 22000
$ws.Range("M23").Value = 22000
$ws.Range("N23").Value = $ws.Range("N24").Value2
$ws.Range("O23").Value = $ws.Range("O24").Value2
$ws.Range("P23").Value = 880
$ws.Range("Q23").Value = $ws.Range("Q24").Value2
$ws.Range("R23").Value = $ws.Range("R24").Value2
